$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final roster: "Nick Smith Jr." was removed and "Ziaire Williams" was added,
# with the remaining players re-sorted into a new order.
$data = @(
  @("Devin Booker","PG,SG","Phoenix Suns"),
  @("Max Christie","SG,SF","Dallas Mavericks"),
  @("Jimmy Butler III","SF,PF","Golden State Warriors"),
  @("LeBron James","SF,PF","Los Angeles Lakers"),
  @("Trae Young","PG","Atlanta Hawks"),
  @("Jalen Williams","SG,SF,PF,C","Oklahoma City Thunder"),
  @("Ziaire Williams","SG,SF","Brooklyn Nets"),
  @("Immanuel Quickley","PG,SG","Toronto Raptors"),
  @("Jalen Brunson","PG","New York Knicks"),
  @("Jabari Smith Jr.","PF,C","Houston Rockets"),
  @("Myles Turner","C","Indiana Pacers"),
  @("Walker Kessler","C","Utah Jazz"),
  @("Kawhi Leonard","SG,SF,PF","LA Clippers"),
  @("Desmond Bane","SG,SF","Memphis Grizzlies"),
  @("Nicolas Claxton","C","Brooklyn Nets"),
  @("Brandon Ingram","SG,SF,PF","Toronto Raptors"),
  @("Norman Powell","SG,SF","LA Clippers"),
  @("Jeremy Sochan","SF,PF","San Antonio Spurs")
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $data[$i][0]
  $ws.Cells.Item($row, 2).Value = $data[$i][1]
  $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
